$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 0
$co = $ws.ChartObjects(1)
$co.TotallyBogusMethodXYZ()
Write-Host "after bogus call"
